{"js": "// Fill in the auto-generated payroll report placeholders with their\n// computed values.\nconst replacements = {\n  countofemployees: \"79\",\n  countofcitizen: \"14\",\n  countofpr: \"25\",\n  countofforeigner: \"40\",\n};\n\nfor (const [placeholder, value] of Object.entries(replacements)) {\n  const found = context.document.body.search(placeholder, { matchCase: true });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(value, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Fill in the auto-generated payroll report placeholders with their\n# computed values (auto fill in word doc).\n$d = $word.ActiveDocument\n\n$replacements = @{\n    \"countofemployees\" = \"79\"\n    \"countofcitizen\"   = \"14\"\n    \"countofpr\"        = \"25\"\n    \"countofforeigner\" = \"40\"\n}\n\nforeach ($key in $replacements.Keys) {\n    $newValue = $replacements[$key]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    # FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    # MatchAllWordForms, Forward, Wrap(wdFindContinue=1), Format, ReplaceWith,\n    # Replace(wdReplaceAll=2)\n    $find.Execute($key, $true, $false, $false, $false, $false, $true, 1, $false, $newValue, 2)\n}\n"}
